$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 150
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 278.2
$ws.Range("I38").Value = 278.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 834.5999999999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -462.5999999999999
$ws.Range("N38").ClearContents()
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20924
$ws.Range("H62").Value = 7000
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 7000
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240
$ws.Range("H80").Value = 805.2
$ws.Range("I80").Value = 717
$ws.Range("K80").Value = 2151
$ws.Range("M80").Value = -1153
$ws.Range("H83").Value = 805.2
$ws.Range("I83").Value = 717
$ws.Range("K83").Value = 6453
$ws.Range("M83").Value = -1461
$ws.Range("H113").Value = 7103
$ws.Range("I113").Value = 5482.25
$ws.Range("K113").Value = 5482.25
$ws.Range("M113").Value = -2228.25
$ws.Range("H135").Value = 999.1667
$ws.Range("I135").Value = 973.75
$ws.Range("K135").Value = 8763.75
$ws.Range("M135").Value = -6228.75
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 2400
$ws.Range("J29").Value = 2400
$ws.Range("L29").Value = 2400
$ws.Range("N29").Value = -3016
$ws.Range("H32").Value = 2252.3333
$ws.Range("I32").Value = 1915.4783
$ws.Range("K32").Value = 1915.4783
$ws.Range("M32").Value = -1628.4783
$ws.Range("H61").Value = 6500
$ws.Range("I61").Value = 6500
$ws.Range("K61").Value = 6500
$ws.Range("M61").Value = -6288
$ws.Range("H92").Value = 9995
$ws.Range("J92").Value = 9995
$ws.Range("L92").Value = 9995
$ws.Range("N92").Value = -14987
$ws.Range("H132").Value = 954
$ws.Range("I132").Value = 954
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2862
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -332
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 6500
$ws.Range("I136").Value = 6500
$ws.Range("K136").Value = 19500
$ws.Range("M136").Value = -16950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 25000
$ws.Range("J9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("N9").Value = -25336
$ws.Range("H134").Value = 11907.6
$ws.Range("I134").Value = 11907.6
$ws.Range("K134").Value = 35722.8
$ws.Range("M134").Value = -33187.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 625.5625
$ws.Range("I107").Value = 580.75
$ws.Range("J107").Value = 760
$ws.Range("K107").Value = 580.75
$ws.Range("L107").Value = 760
$ws.Range("M107").Value = 1339.25
$ws.Range("N107").Value = -4600
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
$ws.Range("H121").Value = 226326
$ws.Range("J121").Value = 226326
$ws.Range("L121").Value = 226326
$ws.Range("N121").Value = -228946
$ws.Range("H122").Value = 1047.7059
$ws.Range("I122").Value = 1137.6666
$ws.Range("K122").Value = 3412.9998
$ws.Range("M122").Value = -962.9998000000001
$ws.Range("H134").Value = 1541.1
$ws.Range("I134").Value = 1451.375
$ws.Range("K134").Value = 4354.125
$ws.Range("M134").Value = -1819.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 577
$ws.Range("I5").Value = 577
$ws.Range("K5").Value = 1731
$ws.Range("M5").Value = -1619
$ws.Range("H37").Value = 99900
$ws.Range("J37").Value = 99900
$ws.Range("L37").Value = 299700
$ws.Range("N37").Value = -299924
$ws.Range("H131").Value = 1991.6666
$ws.Range("I131").Value = 1991.6666
$ws.Range("K131").Value = 5974.9998
$ws.Range("M131").Value = -934.9997999999996
$ws.Range("H134").Value = 200355.6
$ws.Range("I134").Value = 200355.6
$ws.Range("K134").Value = 601066.8
$ws.Range("M134").Value = -595996.8
$ws.Range("H135").Value = 577
$ws.Range("I135").Value = 577
$ws.Range("K135").Value = 5193
$ws.Range("M135").Value = -2658
$ws.Range("H136").Value = 3075
$ws.Range("I136").Value = 3075
$ws.Range("K136").Value = 9225
$ws.Range("M136").Value = -4125
$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("M137").Value = -900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 41669388
$ws.Range("I70").Value = 41669388
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 41669388
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -41669118
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 41669388
$ws.Range("I73").Value = 41669388
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 41669388
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -41668452
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 3562.5
$ws.Range("I80").Value = 3416.6667
$ws.Range("K80").Value = 3416.6667
$ws.Range("M80").Value = -2418.6667
$ws.Range("H83").Value = 3562.5
$ws.Range("I83").Value = 3416.6667
$ws.Range("K83").Value = 17083.3335
$ws.Range("M83").Value = -12091.3335
$ws.Range("H102").Value = 1002.625
$ws.Range("I102").Value = 860.1429000000001
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 860.1429000000001
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 761.8570999999999
$ws.Range("N102").Value = -5244
$ws.Range("H126").Value = 9672.546
$ws.Range("J126").Value = 4000
$ws.Range("L126").Value = 12000
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18966.928
$ws.Range("I7").Value = 15904.2
$ws.Range("K7").Value = 15904.2
$ws.Range("M7").Value = -15792.2
$ws.Range("H32").Value = 10825.5
$ws.Range("I32").Value = 10825.5
$ws.Range("K32").Value = 10825.5
$ws.Range("M32").Value = -10508.5
$ws.Range("H55").Value = 176.66667
$ws.Range("I55").Value = 131.57143
$ws.Range("J55").Value = 239.8
$ws.Range("K55").Value = 131.57143
$ws.Range("L55").Value = 239.8
$ws.Range("M55").Value = 41.42857000000001
$ws.Range("N55").Value = -585.8
$ws.Range("H82").Value = 270
$ws.Range("J82").Value = 250
$ws.Range("L82").Value = 250
$ws.Range("N82").Value = -972
$ws.Range("H85").Value = 270
$ws.Range("J85").Value = 250
$ws.Range("L85").Value = 250
$ws.Range("N85").Value = -2746
$ws.Range("H126").Value = 18966.928
$ws.Range("I126").Value = 15904.2
$ws.Range("K126").Value = 47712.60000000001
$ws.Range("M126").Value = -45242.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 970.46155
$ws.Range("I100").Value = 757
$ws.Range("J100").Value = 1450.75
$ws.Range("K100").Value = 1514
$ws.Range("L100").Value = 2901.5
$ws.Range("M100").Value = -973
$ws.Range("N100").Value = -3983.5
$ws.Range("H112").Value = 19093.5
$ws.Range("J112").Value = 19093.5
$ws.Range("L112").Value = 19093.5
$ws.Range("N112").Value = -22047.5
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 2732.1538
$ws.Range("I132").Value = 2543.1667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7629.500100000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5099.500100000001
$ws.Range("N132").Value = -20060
